$d = $word.ActiveDocument

# 1. Replace the paragraph text "isdifjsdkla" -> "It is changed."
$d.Content.Find.Execute("isdifjsdkla", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "It is changed.", 2)

# 2. Append a brand-new, empty paragraph right after it (before the sectPr),
#    carrying the same run-properties (en-US language) as the existing
#    paragraph mark, but with no run of its own.
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$paragraphXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body>' + `
    '</w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($paragraphXml)
